$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by the Overview summary cells and by the Status
#    column on each per-language detail sheet.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$ws1.Range("E2").Value2 = $newStatus
$ws1.Range("F2").Value2 = $newStatus
$ws1.Range("E3").Value2 = $newStatus
$ws1.Range("F3").Value2 = $newStatus

$ws2.Range("C2").Value2 = $newStatus
$ws2.Range("C3").Value2 = $newStatus

$ws3.Range("C2").Value2 = $newStatus
$ws3.Range("C3").Value2 = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both handed-back documents.
# ---------------------------------------------------------------------------
$ws2.Range("I2").Value2 = "00445e6c-6b79-4512-aab5-d0bbbd82c7f4.md"
$ws2.Range("J2").Value2 = "00445e6c-6b79-4512-aab5-d0bbbd82c7f4.efdffa79b79ca2ab551b8a38d72c6df99778b770.zh-cn.xlf"
$ws2.Range("K2").Value2 = "2016-08-27 10:47:42"

$ws2.Range("I3").Value2 = "0ba3c325-b4ba-4a0f-bbfa-54202c9eb50a.md"
$ws2.Range("J3").Value2 = "0ba3c325-b4ba-4a0f-bbfa-54202c9eb50a.27d8ea604e44cee68fc1b1363eb379a21877326c.zh-cn.xlf"
$ws2.Range("K3").Value2 = "2016-08-27 10:47:42"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same fields, different handback timestamp.
# ---------------------------------------------------------------------------
$ws3.Range("I2").Value2 = "00445e6c-6b79-4512-aab5-d0bbbd82c7f4.md"
$ws3.Range("J2").Value2 = "00445e6c-6b79-4512-aab5-d0bbbd82c7f4.efdffa79b79ca2ab551b8a38d72c6df99778b770.de-de.xlf"
$ws3.Range("K2").Value2 = "2016-08-27 10:47:49"

$ws3.Range("I3").Value2 = "0ba3c325-b4ba-4a0f-bbfa-54202c9eb50a.md"
$ws3.Range("J3").Value2 = "0ba3c325-b4ba-4a0f-bbfa-54202c9eb50a.27d8ea604e44cee68fc1b1363eb379a21877326c.de-de.xlf"
$ws3.Range("K3").Value2 = "2016-08-27 10:47:49"

# ---------------------------------------------------------------------------
# 4. Rebuild the hyperlinks on the detail sheets so that the newly
#    populated "Latest Target File" cells (I2/I3) link to the same source
#    document pages as A2/A3, preserving relationship ordering
#    (A2, I2, A3, I3).
# ---------------------------------------------------------------------------
$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1308ade67f2860791d17a79989cf4f22481a7356/e2e/00445e6c-6b79-4512-aab5-d0bbbd82c7f4.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1308ade67f2860791d17a79989cf4f22481a7356/e2e/0ba3c325-b4ba-4a0f-bbfa-54202c9eb50a.md"

foreach ($ws in @($ws2, $ws3)) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $url1, $null, $null, "00445e6c-6b79-4512-aab5-d0bbbd82c7f4.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $url1, $null, $null, "00445e6c-6b79-4512-aab5-d0bbbd82c7f4.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $url2, $null, $null, "0ba3c325-b4ba-4a0f-bbfa-54202c9eb50a.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $url2, $null, $null, "0ba3c325-b4ba-4a0f-bbfa-54202c9eb50a.md")

    # Give the new "Latest Target File" hyperlinks the same visual style as
    # the existing "Source File Name" hyperlinks in column A.
    $ws.Range("I2:I3").Style = "HyperLink"
}

# ---------------------------------------------------------------------------
# 5. Column width adjustments on the two detail sheets (Status column grew
#    wider because of the new status text; the two newly populated columns
#    grew to their maximum auto-fit width of 40 characters).
# ---------------------------------------------------------------------------
foreach ($ws in @($ws2, $ws3)) {
    $ws.Columns.Item(3).ColumnWidth = 29.17
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}

$wb.Save()
